$d = $word.ActiveDocument

# Locate the bracketing paragraphs for the bullet-list block that needs to be
# rewritten, by matching on text rather than a hard-coded index, so the
# script is resilient to any paragraphs before this block.
#   - $startAnchorIdx : "Regression models with categorical covariates"
#                        (last untouched paragraph before the block)
#   - $lastOldIdx      : "Peer assessment: ..." (last paragraph of the block
#                        as it exists today)
$count = $d.Paragraphs.Count
$startAnchorIdx = -1
$lastOldIdx = -1
for ($i = 1; $i -le $count; $i++) {
    $t = $d.Paragraphs($i).Range.Text
    if ($t -like "Regression models*") {
        $startAnchorIdx = $i
    }
    if ($t -like "Peer assessment*") {
        $lastOldIdx = $i
    }
}

# Delete the whole old block (6 paragraphs: "Re-do notes for Github" through
# "Peer assessment ... marking"), including each paragraph's own paragraph
# mark, by spanning from the start of the first old paragraph to the start
# of the paragraph that follows the block. Deleting this way (rather than
# clearing each paragraph's Range.Text individually) also drops the stray
# w:proofErr spell/grammar-check markers that used to split some of these
# paragraphs into multiple runs.
$firstOldPara = $d.Paragraphs($startAnchorIdx + 1)
$afterOldPara = $d.Paragraphs($lastOldIdx + 1)
$oldRange = $d.Range($firstOldPara.Range.Start, $afterOldPara.Range.Start)
$oldRange.Delete()

# The updated / expanded set of bullet items that replace the old block.
$items = @(
    "Re-do notes for Github",
    "Add latex external resources",
    "Add Quarto and Git as assessed material on quiz 2",
    "Increase the marks for version control assessment (in groups projects)",
    "Make quarto slides mandatory for peer review assessment.",
    "More marks to git.",
    "Peer assessment: make sure they really use tidymodels and put it into the marking.",
    "Word count limit rather than page limit for the group project.",
    "(potentially) Add classification & Random Forests/Bagging material.",
    "Teach ANOVA’s?",
    "Show raw marks for each assignment, e.g. 16 points scale for group projects."
)

# Re-insert one paragraph per item, anchored right after the untouched
# "Regression models ..." paragraph, so each new paragraph inherits its
# ListParagraph style + list numbering (numPr) via InsertParagraphAfter.
$anchorIndex = $startAnchorIdx
foreach ($item in $items) {
    $anchor = $d.Paragraphs($anchorIndex)
    $anchor.Range.InsertParagraphAfter()
    $anchorIndex = $anchorIndex + 1
    $d.Paragraphs($anchorIndex).Range.Text = $item
}
